# Natmi following Dr Hou advice
#
# The sending/target cluster set grows from {ECs, sCs} to {ECs, FAPs, sCs},
# so the ligand-receptor pair table (F11r -> Jam2) now has 3 x 3 = 9 rows
# (rows 2-10) instead of 2 x 3 = 6 (rows 2-7), and every numeric column is
# recomputed for the new cluster composition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns, in order: A..T
#  A Sending cluster                                  K Receptor-expressing cells
#  B Ligand symbol                                    L Receptor detection rate
#  C Receptor symbol                                  M Receptor average expression value
#  D Target cluster                                   N Receptor total expression value
#  E Ligand-expressing cells                           O Receptor derived specificity (avg)
#  F Ligand detection rate                             P Receptor derived specificity (total)
#  G Ligand average expression value                   Q Edge average expression weight
#  H Ligand total expression value                     R Edge total expression weight
#  I Ligand derived specificity (avg)                   S Edge average expression derived specificity
#  J Ligand derived specificity (total)                 T Edge total expression derived specificity
$rows = @(
  @("ECs","F11r","Jam2","ECs",3,1,20.10858766666667,60.325763,0.9234633461941997,0.9234633461941997,3,1,33.24999533333333,99.74998599999999,0.5673360890306117,0.5673360890306117,668.6104460765908,6017.494014689318,0.5239140831929391,0.5239140831929391),
  @("ECs","F11r","Jam2","FAPs",3,1,20.10858766666667,60.325763,0.9234633461941997,0.9234633461941997,3,1,23.30243966666667,69.907319,0.3976035140102714,0.3976035140102714,468.5791508843775,4217.212357959397,0.3671722715064976,0.3671722715064976),
  @("ECs","F11r","Jam2","sCs",3,1,20.10858766666667,60.325763,0.9234633461941997,0.9234633461941997,3,1,2.054792666666667,6.164378,0.03506039695911681,0.03506039695911681,41.31897847449045,371.870806270414,0.03237699149476295,0.03237699149476295),
  @("FAPs","F11r","Jam2","ECs",2,0.6666666666666666,0.2128156666666667,0.638447,0.009773310334883756,0.009773310334883756,3,1,33.24999533333333,99.74998599999999,0.5673360890306117,0.5673360890306117,7.076119923526888,63.68507931174199,0.005544751662275409,0.005544751662275409),
  @("FAPs","F11r","Jam2","FAPs",2,0.6666666666666666,0.2128156666666667,0.638447,0.009773310334883756,0.009773310334883756,3,1,23.30243966666667,69.907319,0.3976035140102714,0.3976035140102714,4.959124232621445,44.632118093593,0.003885902532662684,0.003885902532662684),
  @("FAPs","F11r","Jam2","sCs",2,0.6666666666666666,0.2128156666666667,0.638447,0.009773310334883756,0.009773310334883756,3,1,2.054792666666667,6.164378,0.03506039695911681,0.03506039695911681,0.4372920712184444,3.935628640966,0.0003426561399456634,0.0003426561399456633),
  @("sCs","F11r","Jam2","ECs",3,1,1.453784333333333,4.361352999999999,0.06676334347091657,0.06676334347091657,3,1,33.24999533333333,99.74998599999999,0.5673360890306117,0.5673360890306117,48.33832229900644,435.0449006910579,0.03787725417539724,0.03787725417539724),
  @("sCs","F11r","Jam2","FAPs",3,1,1.453784333333333,4.361352999999999,0.06676334347091657,0.06676334347091657,3,1,23.30243966666667,69.907319,0.3976035140102714,0.3976035140102714,33.87672171584522,304.8904954426069,0.02654533997111114,0.02654533997111114),
  @("sCs","F11r","Jam2","sCs",3,1,1.453784333333333,4.361352999999999,0.06676334347091657,0.06676334347091657,3,1,2.054792666666667,6.164378,0.03506039695911681,0.03506039695911681,2.987225387048222,26.885028483434,0.002340749324408195,0.002340749324408194)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $rowNum = $i + 2
    for ($c = 0; $c -lt 20; $c++) {
        $val = $r[$c]
        $ws.Cells.Item($rowNum, $c + 1).Value = $val
    }
}
